# Apply French translation + two additional risk rows to the smartfridge risks sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Risque"
$ws.Range("B1").Value = "Probabilité"
$ws.Range("C1").Value = "Impact"
$ws.Range("D1").Value = "Mitigation"

# Row 2
$ws.Range("A2").Value = "Mauvaise reconnaissance produit par l'IA"
$ws.Range("B2").Value = "Moyenne"
$ws.Range("C2").Value = "Élevé (perte confiance utilisateur)"
$ws.Range("D2").Value = "Phase apprentissage supervisée + correction manuelle utilisateur"

# Row 3
$ws.Range("A3").Value = "Dépendance à un seul Drive (ex : Leclerc)"
$ws.Range("B3").Value = "Moyenne"
$ws.Range("C3").Value = "Élevé (killer feature KO)"
$ws.Range("D3").Value = "Prévoir intégration multi-retailers dès la conception (Carrefour, Intermarché)"

# Row 4
$ws.Range("A4").Value = "Rejet utilisateur (caméra perçue comme intrusive)"
$ws.Range("B4").Value = "Élevée"
$ws.Range("C4").Value = "Élevé (frein adoption marché)"
$ws.Range("D4").Value = "Communication claire : aucune image stockée, IA locale uniquement"

# Row 5
$ws.Range("A5").Value = "Problèmes RGPD / conformité données"
$ws.Range("B5").Value = "Élevée"
$ws.Range("C5").Value = "Élevé (risque légal)"
$ws.Range("D5").Value = "Edge AI : aucune image brute envoyée au cloud, chiffrement, Privacy Officer"

# Row 6
$ws.Range("A6").Value = "Perte de connexion Wi-Fi entre frigo et cloud"
$ws.Range("B6").Value = "Élevée"
$ws.Range("C6").Value = "Moyen (perte synchro panier auto)"
$ws.Range("D6").Value = "Cache local + re-sync dès que Wi-Fi revient"

# Row 7 (new)
$ws.Range("A7").Value = "Dépassement budget matériel (coûts IoT)"
$ws.Range("B7").Value = "Moyenne"
$ws.Range("C7").Value = "Moyen (retarde le prototype)"
$ws.Range("D7").Value = "Buffer financier 10% + achat alternatif composants"

# Row 8 (new)
$ws.Range("A8").Value = "Délai ou blocage négociation accès catalogue Drive"
$ws.Range("B8").Value = "Moyenne"
$ws.Range("C8").Value = "Élevé (bloque le go-to-market)"
$ws.Range("D8").Value = "Business Dev dédié très tôt dans le projet"
